$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 380.08334
$ws.Range("J96").Value = 98.5
$ws.Range("L96").Value = 295.5
$ws.Range("N96").Value = -3041.5
$ws.Range("H98").Value = 611.37036
$ws.Range("I98").Value = 611.37036
$ws.Range("K98").Value = 611.37036
$ws.Range("M98").Value = 886.62964
$ws.Range("H122").Value = 611.37036
$ws.Range("I122").Value = 611.37036
$ws.Range("K122").Value = 1834.11108
$ws.Range("M122").Value = 615.8889199999999
$ws.Range("H132").Value = 4173.3057
$ws.Range("I132").Value = 918.5925999999999
$ws.Range("J132").Value = 13937.444
$ws.Range("K132").Value = 2755.7778
$ws.Range("L132").Value = 41812.33199999999
$ws.Range("M132").Value = -225.7777999999998
$ws.Range("N132").Value = -46872.33199999999
$ws.Range("H135").Value = 1201.4117
$ws.Range("I135").Value = 1245.25
$ws.Range("K135").Value = 11207.25
$ws.Range("M135").Value = -8672.25
$ws.Range("H137").Value = 1626.7142
$ws.Range("I137").Value = 1221.25
$ws.Range("K137").Value = 3663.75
$ws.Range("M137").Value = -1113.75
$ws.Range("H138").Value = 2654.52
$ws.Range("I138").Value = 2591
$ws.Range("K138").Value = 7773
$ws.Range("M138").Value = -2633
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1447.5714
$ws.Range("I32").Value = 1447.5714
$ws.Range("K32").Value = 1447.5714
$ws.Range("M32").Value = -1160.5714
$ws.Range("H122").Value = 2274
$ws.Range("I122").Value = 1712.6666
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 5137.9998
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -2687.9998
$ws.Range("N122").Value = -19300
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13999
$ws.Range("I26").Value = 13999
$ws.Range("K26").Value = 13999
$ws.Range("M26").Value = -13707
$ws.Range("H35").Value = 31028.334
$ws.Range("J35").Value = 32074
$ws.Range("L35").Value = 32074
$ws.Range("N35").Value = -32694
$ws.Range("H96").Value = 30000
$ws.Range("I96").Value = 30000
$ws.Range("K96").Value = 30000
$ws.Range("M96").Value = -27254
$ws.Range("H134").Value = 3600.7112
$ws.Range("I134").Value = 3600.7112
$ws.Range("K134").Value = 10802.1336
$ws.Range("M134").Value = -8267.133600000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 25587
$ws.Range("I10").Value = 5050
$ws.Range("J10").Value = 66661
$ws.Range("K10").Value = 5050
$ws.Range("L10").Value = 66661
$ws.Range("M10").Value = -4911
$ws.Range("N10").Value = -66939
$ws.Range("H22").Value = 261
$ws.Range("I22").Value = 268.1
$ws.Range("K22").Value = 268.1
$ws.Range("M22").Value = 81.89999999999998
$ws.Range("H31").Value = 4425.5293
$ws.Range("I31").Value = 2295.5386
$ws.Range("K31").Value = 2295.5386
$ws.Range("M31").Value = -2000.5386
$ws.Range("H34").Value = 4425.5293
$ws.Range("I34").Value = 2295.5386
$ws.Range("K34").Value = 2295.5386
$ws.Range("M34").Value = -2093.5386
$ws.Range("H107").Value = 566.65955
$ws.Range("I107").Value = 530.125
$ws.Range("J107").Value = 644.6
$ws.Range("K107").Value = 530.125
$ws.Range("L107").Value = 644.6
$ws.Range("M107").Value = 1389.875
$ws.Range("N107").Value = -4484.6
$ws.Range("H134").Value = 5947.706
$ws.Range("I134").Value = 4740.7334
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 14222.2002
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -11687.2002
$ws.Range("N134").Value = -50070
$ws.Range("H141").Value = 27882.166
$ws.Range("I141").Value = 25458.6
$ws.Range("K141").Value = 25458.6
$ws.Range("M141").Value = -20278.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1111849.8
$ws.Range("J117").Value = 1250718.5
$ws.Range("L117").Value = 3752155.5
$ws.Range("N117").Value = -3759039.5
$ws.Range("H121").Value = 7214817
$ws.Range("J121").Value = 14429187
$ws.Range("L121").Value = 43287561
$ws.Range("N121").Value = -43290181
$ws.Range("H126").Value = 6331.6
$ws.Range("I126").Value = 4164.75
$ws.Range("K126").Value = 12494.25
$ws.Range("M126").Value = -7554.25
$ws.Range("H128").Value = 422994.22
$ws.Range("I128").Value = 422994.22
$ws.Range("K128").Value = 1268982.66
$ws.Range("M128").Value = -1264002.66
$ws.Range("H129").Value = 653601.9399999999
$ws.Range("I129").Value = 112574.664
$ws.Range("J129").Value = 1262257.6
$ws.Range("K129").Value = 337723.992
$ws.Range("L129").Value = 3786772.8
$ws.Range("M129").Value = -332723.992
$ws.Range("N129").Value = -3796772.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 49600
$ws.Range("I62").Value = 42900
$ws.Range("J62").Value = 63000
$ws.Range("K62").Value = 42900
$ws.Range("L62").Value = 63000
$ws.Range("M62").Value = -42214
$ws.Range("N62").Value = -64372
$ws.Range("H65").Value = 49600
$ws.Range("I65").Value = 42900
$ws.Range("J65").Value = 63000
$ws.Range("K65").Value = 128700
$ws.Range("L65").Value = 189000
$ws.Range("M65").Value = -125268
$ws.Range("N65").Value = -195864
$ws.Range("H102").Value = 2617.7368
$ws.Range("I102").Value = 2072.7646
$ws.Range("K102").Value = 2072.7646
$ws.Range("M102").Value = -450.7646
$ws.Range("H134").Value = 53998.4
$ws.Range("J134").Value = 53998.4
$ws.Range("L134").Value = 161995.2
$ws.Range("N134").Value = -167065.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 483
$ws.Range("I40").Value = 483
$ws.Range("K40").Value = 483
$ws.Range("M40").Value = -347
$ws.Range("H55").Value = 549.3684
$ws.Range("I55").Value = 628.25
$ws.Range("J55").Value = 414.14285
$ws.Range("K55").Value = 628.25
$ws.Range("L55").Value = 414.14285
$ws.Range("M55").Value = -455.25
$ws.Range("N55").Value = -760.14285
$ws.Range("H74").Value = 48666.668
$ws.Range("I74").Value = 39000
$ws.Range("K74").Value = 39000
$ws.Range("M74").Value = -38002
$ws.Range("H77").Value = 48666.668
$ws.Range("I77").Value = 39000
$ws.Range("K77").Value = 117000
$ws.Range("M77").Value = -112008
$ws.Range("H100").Value = 4961.154
$ws.Range("I100").Value = 3686.875
$ws.Range("K100").Value = 3686.875
$ws.Range("M100").Value = -3145.875
$ws.Range("H122").Value = 4254.5557
$ws.Range("J122").Value = 4599.5
$ws.Range("L122").Value = 13798.5
$ws.Range("N122").Value = -18698.5
$ws.Range("H132").Value = 9865.333000000001
$ws.Range("I132").Value = 9930.936
$ws.Range("K132").Value = 29792.808
$ws.Range("M132").Value = -27262.808
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1245614.9
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1245614.9
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1245614.9
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1245838.9
$ws.Range("H81").Value = 19998.5
$ws.Range("I81").Value = 19995
$ws.Range("K81").Value = 39990
$ws.Range("M81").Value = -38929
$ws.Range("H84").Value = 19998.5
$ws.Range("I84").Value = 19995
$ws.Range("K84").Value = 199950
$ws.Range("M84").Value = -194646
$ws.Range("H93").Value = 10000
$ws.Range("J93").Value = 10000
$ws.Range("L93").Value = 10000
$ws.Range("N93").Value = -14992
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H126").Value = 3177.75
$ws.Range("I126").Value = 3177.75
$ws.Range("K126").Value = 9533.25
$ws.Range("M126").Value = -7063.25
$ws.Range("H132").Value = 3970.0908
$ws.Range("I132").Value = 3602.7307
$ws.Range("J132").Value = 5334.5713
$ws.Range("K132").Value = 10808.1921
$ws.Range("L132").Value = 16003.7139
$ws.Range("M132").Value = -8278.1921
$ws.Range("N132").Value = -21063.7139
$ws.Range("H136").Value = 3055.4546
$ws.Range("I136").Value = 2358.1724
$ws.Range("J136").Value = 8110.75
$ws.Range("K136").Value = 7074.5172
$ws.Range("L136").Value = 24332.25
$ws.Range("M136").Value = -4524.5172
$ws.Range("N136").Value = -29432.25
